$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "Identifier" -> "Stock Code"
$ws.Range("A1").Value = "Stock Code"

# New column H with header + values, matching header style of existing header row
$ws.Range("H1").Value = "Column Name"
$ws.Range("H1").Style = $ws.Range("G1").Style

$ws.Range("H2").Value = "Open"
$ws.Range("H3").Value = "Open"
$ws.Range("H4").Value = "Close"
$ws.Range("H5").Value = "Close"
$ws.Range("H6").Value = "Close"
$ws.Range("H7").Value = "Volume"
$ws.Range("H8").Value = "Market Cap"
$ws.Range("H9").Value = "Market Cap"
$ws.Range("H10").Value = "Market Cap"
$ws.Range("H11").Value = "Market Cap"

# Category column: ACCEPTABLE -> Warning for rows 8-11 (Market_Cap_check rows)
$ws.Range("C8").Value = "Warning"
$ws.Range("C9").Value = "Warning"
$ws.Range("C10").Value = "Warning"
$ws.Range("C11").Value = "Warning"
